# hed3_tags_single_sheet.xlsx edit
# The HED tags column (D) entry for the "PerturbLeft" event (row 2) was
# replaced with a new, shorter placeholder tag string.  Everything else on
# the sheet (the other rows, headers, etc.) stays the same - only the text
# shown in D2 changes, which also means the row auto-shrinks now that the
# text is much shorter, and the active cell selection moved to D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the HED tag text for the "PerturbLeft" row.
$ws.Range("D2").Value = "Attribute/Sensory/Bisual"

# The row no longer needs to be tall enough for the long HED string, so its
# height shrinks down close to the sheet's default row height.
$ws.Rows.Item(2).RowHeight = 14.9

# Move the active selection to D3 (matches the saved sheet view state).
$ws.Range("D3").Select() | Out-Null
